# Updated cryptos list on Sat Jul 15 18:40:03 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay text-typed (matches the source inlineStr cells)
    # instead of letting Excel auto-convert a numeric-looking string (e.g.
    # "0.9991") into a real number.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
}

$ws.Range("D2").Value = "30.331.53"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "1.936.03"
$ws.Range("E3").Value = "  -0.50%  "
Set-TextValue $ws.Range("D4") "0.9991"
$ws.Range("E4").Value = "  -1.25%  "
Set-TextValue $ws.Range("D5") "251.60"
$ws.Range("E5").Value = "  +0.92%  "
Set-TextValue $ws.Range("D6") "0.7184"
$ws.Range("E6").Value = "  +4.01%  "
Set-TextValue $ws.Range("D7") "0.9995"
$ws.Range("E7").Value = "  -1.26%  "
Set-TextValue $ws.Range("D8") "0.3309"
$ws.Range("E8").Value = "  +1.80%  "
Set-TextValue $ws.Range("D9") "27.97"
$ws.Range("E9").Value = "  +5.64%  "
Set-TextValue $ws.Range("D10") "0.07293"
$ws.Range("E10").Value = "  +7.23%  "
Set-TextValue $ws.Range("D11") "0.8064"
$ws.Range("E11").Value = "  +1.03%  "
Set-TextValue $ws.Range("D12") "0.08106"
$ws.Range("E12").Value = "  +1.18%  "
$ws.Range("D13").Value = "1.932.51"
$ws.Range("E13").Value = "  -0.86%  "
Set-TextValue $ws.Range("D14") "5.491"
$ws.Range("E14").Value = "  +1.38%  "
Set-TextValue $ws.Range("D15") "94.90"
$ws.Range("E15").Value = "  +0.54%  "
Set-TextValue $ws.Range("D16") "15.13"
$ws.Range("E16").Value = "  +3.92%  "
$ws.Range("D17").Value = "30.311.60"
$ws.Range("E17").Value = "  -0.87%  "
Set-TextValue $ws.Range("D18") "253.73"
$ws.Range("E18").Value = "  -3.43%  "
Set-TextValue $ws.Range("D19") "0.000008214"
$ws.Range("E19").Value = "  +4.56%  "
Set-TextValue $ws.Range("D20") "5.831"
$ws.Range("E20").Value = "  -1.22%  "
$ws.Range("D21").Value = "2.185.16"
$ws.Range("E21").Value = "  -1.59%  "
Set-TextValue $ws.Range("D22") "0.9996"
$ws.Range("E22").Value = "  -1.08%  "
Set-TextValue $ws.Range("D23") "0.9986"
$ws.Range("E23").Value = "  -1.31%  "
Set-TextValue $ws.Range("D24") "6.995"
$ws.Range("E24").Value = "  +1.84%  "
Set-TextValue $ws.Range("D25") "9.774"
$ws.Range("E25").Value = "  +0.78%  "
Set-TextValue $ws.Range("D26") "164.95"
$ws.Range("E26").Value = "  +3.87%  "
Set-TextValue $ws.Range("D27") "2.365"
$ws.Range("E27").Value = "  +4.06%  "
Set-TextValue $ws.Range("D28") "19.38"
$ws.Range("E28").Value = "  +2.65%  "
Set-TextValue $ws.Range("D29") "0.1309"
$ws.Range("E29").Value = "  +1.22%  "
Set-TextValue $ws.Range("D30") "1.352"
$ws.Range("E30").Value = "  -2.21%  "
Set-TextValue $ws.Range("D31") "1.540"
$ws.Range("E31").Value = "  -1.77%  "
Set-TextValue $ws.Range("D32") "4.442"
$ws.Range("E32").Value = "  +0.27%  "
Set-TextValue $ws.Range("D33") "4.201"
$ws.Range("E33").Value = "  -1.07%  "
Set-TextValue $ws.Range("D34") "0.05239"
$ws.Range("E34").Value = "  +2.40%  "
Set-TextValue $ws.Range("D35") "1.271"
$ws.Range("E35").Value = "  +5.58%  "
Set-TextValue $ws.Range("D36") "0.7504"
$ws.Range("E36").Value = "  -0.14%  "
Set-TextValue $ws.Range("D37") "2.772"
$ws.Range("E37").Value = "  +1.02%  "
Set-TextValue $ws.Range("D38") "0.01978"
$ws.Range("E38").Value = "  +1.67%  "
Set-TextValue $ws.Range("D39") "2.812"
$ws.Range("E39").Value = "  -0.18%  "
Set-TextValue $ws.Range("D40") "79.18"
$ws.Range("E40").Value = "  -1.41%  "
Set-TextValue $ws.Range("D41") "6.454"
$ws.Range("E41").Value = "  -2.26%  "
Set-TextValue $ws.Range("D42") "0.4553"
$ws.Range("E42").Value = "  +2.50%  "
Set-TextValue $ws.Range("D43") "2.037"
$ws.Range("E43").Value = "  -0.81%  "
Set-TextValue $ws.Range("D44") "0.8439"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("E45").Value = "  -1.14%  "
Set-TextValue $ws.Range("D46") "101.84"
$ws.Range("E46").Value = "  +0.08%  "
Set-TextValue $ws.Range("D47") "9.778"
$ws.Range("E47").Value = "  +0.41%  "
Set-TextValue $ws.Range("D48") "7.474"
$ws.Range("E48").Value = "  +2.06%  "
Set-TextValue $ws.Range("D49") "36.94"
$ws.Range("E49").Value = "  +2.06%  "
Set-TextValue $ws.Range("D50") "0.4201"
$ws.Range("E50").Value = "  +2.32%  "

# Row 51 is fully replaced: coin name/link/price/volume all change
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D51") "1.501"
$ws.Range("E51").Value = "  +0.06%  "

